$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last used row on the sheet (data rows 2..LastRow, header in row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Insert a new column before GY. This shifts the existing "nom" (GY) and
# "url_produit" (GZ) columns one position to the right (to GZ and HA),
# preserving their header styling/formatting automatically.
$ws.Range("GY1").EntireColumn.Insert()

# New column header: latest snapshot timestamp.
$ws.Range("GY1").Value = "2026-02-06 14:34:03"

# Populate the new price snapshot column by copying the previous/latest
# snapshot column (still GX, since we inserted to its right) for every
# data row. This carries forward each product's most recent known price
# (and keeps rows with no price yet blank), matching the column's
# per-cell formatting exactly.
$ws.Range("GX2:GX$lastRow").Copy($ws.Range("GY2:GY$lastRow"))
